$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sayfa1")

# Row 11: Correctness went from 1.0 -> 0.0, and a Correctness-Reasons note
# ("2.5.1") was added where there was none before.
$ws.Range("F11").Value = 0.0
$ws.Range("G11").Value = "2.5.1"
